# Fixing path variable in Excel: the "CatPose" row label incorrectly read
# "CatPoseVideoPose" (should mirror the "TreePoseVideoPath" label pattern
# used elsewhere in the sheet) - correct it to "CatPoseVideoPath".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "CatPoseVideoPath"

# Reposition the view/selection to match the refactored layout: scrolled so
# row 4 is the top-most visible row, with A7 as the active selected cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("A7").Select()
